# Update biomass/product/substrate/volume training data on worksheet "0"
# per the fixed ODE integration (bug fix in custom ODE function).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("0")

$Cvals = @(0.2037148833171987,0.2435010902287608,0.3284695469044285,0.4661307119061669,0.6677019016349858,0.9281864050651976,1.202536059833841,1.427754454545161,1.593240725795179,1.727499544555397,1.846875227883147,1.956323344284391,2.057737286990446,2.152198756660972,2.240522646756836,2.32334686314293,2.401191819506123,2.474504210661007,2.543680905080736,2.609126124480566)
$Dvals = @(62.24872802469813,59.15949792686844,56.37309103471762,53.85108408086009,51.56274117981751,49.47954851620818,47.56936927436864,45.80002822765291,44.15271568055753,42.61882219470806,41.18927383881763,39.85436000839648,38.60513472872117,37.43365216134531,36.33289991112127,35.29667223726661,34.31946194820772,33.39636910694175,32.52302194988601,31.69551767660233)
$Evals = @(0.1783640750866079,0.5675025526327994,0.821816497379465,0.9323270712251929,0.8815905502993293,0.6837642420872594,0.4324444089877149,0.2565863869787361,0.1799176012720223,0.1475943611260985,0.1287959214725349,0.1152306830680628,0.1045229999830395,0.09574468904220901,0.08835184399556577,0.0820271023945329,0.07657881552199949,0.07186864267836592,0.06777645696811314,0.06409667732964447)
$Fvals = @(4.017967202564075,4.228493518353549,4.439019834143022,4.649546149932496,4.860072465721969,5.070598781511443,5.281125097300917,5.49165141309039,5.702177728879865,5.912704044669339,6.123230360458812,6.333756676248285,6.544282992037759,6.754809307827233,6.965335623616706,7.17586193940618,7.386388255195653,7.596914570985128,7.807440886774601,8.017967202564074)

for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $Cvals[$i]
    $ws.Cells.Item($row, 4).Value = $Dvals[$i]
    $ws.Cells.Item($row, 5).Value = $Evals[$i]
    $ws.Cells.Item($row, 6).Value = $Fvals[$i]
}
